$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-7 with combined tuple-like string values
$ws.Range("A2").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '5/5'])"
$ws.Range("A3").Value = "('Elephant', ['Token Creature — Elephant', '3/3'])"
$ws.Range("A4").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"
$ws.Range("A5").Value = "('Squirrel', ['Token Creature — Squirrel', '1/1'])"
$ws.Range("A6").Value = "('Wurm', ['Token Creature — Wurm', '6/6'])"
$ws.Range("A7").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"

# Delete the now-obsolete rows 8-20
$ws.Range("A8:A20").EntireRow.Delete()
